$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.4912533333333333
$ws.Range("H2").Value = 1.47376
$ws.Range("I2").Value = 0.1305673009957422
$ws.Range("J2").Value = 0.1305673009957422
$ws.Range("P2").Value = 0.9810128591839974
$ws.Range("Q2").Value = 0.00505941808
$ws.Range("R2").Value = 0.04553476272
$ws.Range("S2").Value = 0.1280882012657707
$ws.Range("T2").Value = 0.1280882012657706
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.4912533333333333
$ws.Range("H3").Value = 1.47376
$ws.Range("I3").Value = 0.1305673009957422
$ws.Range("J3").Value = 0.1305673009957422
$ws.Range("Q3").Value = 0.00009792316444444445
$ws.Range("R3").Value = 0.00088130848
$ws.Range("S3").Value = 0.002479099729971547
$ws.Range("T3").Value = 0.002479099729971546
$ws.Range("G4").Value = 0.2760593333333333
$ws.Range("H4").Value = 0.828178
$ws.Range("I4").Value = 0.07337216792697034
$ws.Range("J4").Value = 0.07337216792697034
$ws.Range("P4").Value = 0.9810128591839974
$ws.Range("Q4").Value = 0.002843135074
$ws.Range("R4").Value = 0.025588215666
$ws.Range("S4").Value = 0.07197904024256557
$ws.Range("T4").Value = 0.07197904024256556
$ws.Range("G5").Value = 0.2760593333333333
$ws.Range("H5").Value = 0.828178
$ws.Range("I5").Value = 0.07337216792697034
$ws.Range("J5").Value = 0.07337216792697034
$ws.Range("Q5").Value = 0.00005502782711111111
$ws.Range("R5").Value = 0.000495250444
$ws.Range("S5").Value = 0.001393127684404771
$ws.Range("T5").Value = 0.001393127684404771
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 1.058191
$ws.Range("H6").Value = 3.174573
$ws.Range("I6").Value = 0.2812502907013058
$ws.Range("J6").Value = 0.2812502907013057
$ws.Range("P6").Value = 0.9810128591839974
$ws.Range("Q6").Value = 0.010898309109
$ws.Range("R6").Value = 0.098084781981
$ws.Range("S6").Value = 0.2759101518272185
$ws.Range("T6").Value = 0.2759101518272184
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 1.058191
$ws.Range("H7").Value = 3.174573
$ws.Range("I7").Value = 0.2812502907013058
$ws.Range("J7").Value = 0.2812502907013057
$ws.Range("Q7").Value = 0.0002109327393333334
$ws.Range("R7").Value = 0.001898394654
$ws.Range("S7").Value = 0.005340138874087343
$ws.Range("T7").Value = 0.005340138874087342
$ws.Range("G8").Value = 0.306452
$ws.Range("H8").Value = 0.9193560000000001
$ws.Range("I8").Value = 0.08145005399402996
$ws.Range("J8").Value = 0.08145005399402996
$ws.Range("P8").Value = 0.9810128591839974
$ws.Range("Q8").Value = 0.003156149148
$ws.Range("R8").Value = 0.028405342332
$ws.Range("S8").Value = 0.07990355034937431
$ws.Range("T8").Value = 0.0799035503493743
$ws.Range("G9").Value = 0.306452
$ws.Range("H9").Value = 0.9193560000000001
$ws.Range("I9").Value = 0.08145005399402996
$ws.Range("J9").Value = 0.08145005399402996
$ws.Range("Q9").Value = 0.00006108609866666667
$ws.Range("R9").Value = 0.0005497748880000001
$ws.Range("S9").Value = 0.001546503644655657
$ws.Range("T9").Value = 0.001546503644655657
$ws.Range("G10").Value = 1.630497333333333
$ws.Range("H10").Value = 4.891492
$ws.Range("I10").Value = 0.4333601863819517
$ws.Range("J10").Value = 0.4333601863819517
$ws.Range("P10").Value = 0.9810128591839974
$ws.Range("Q10").Value = 0.016792492036
$ws.Range("R10").Value = 0.151132428324
$ws.Range("S10").Value = 0.4251319154990685
$ws.Range("T10").Value = 0.4251319154990684
$ws.Range("G11").Value = 1.630497333333333
$ws.Range("H11").Value = 4.891492
$ws.Range("I11").Value = 0.4333601863819517
$ws.Range("J11").Value = 0.4333601863819517
$ws.Range("Q11").Value = 0.0003250124684444444
$ws.Range("R11").Value = 0.002925112216
$ws.Range("S11").Value = 0.008228270882883222
$ws.Range("T11").Value = 0.008228270882883222

Write-Output "Applied all cell updates"
